$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TAYDA ORDER")

# Update the Unit Price for the "Shipping & Handling" line (row 16).
$ws.Range("E16").Value = 13.87

# Reflect the cell selection that was active when the file was saved.
$ws.Range("B17").Select()
